$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C:L -> D:M)
$ws.Columns("C:C").Insert()

# New header in C1, matching the formatting of the other header cells
$ws.Range("C1").Value2 = "statut_name"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Map statut_label (column B) -> statut_name (column C) text per row
$map = @{
    "noir"   = "pas de résultat ni de publication"
    "rouge"  = "résultat et / ou publication posté"
    "vert"   = "résultat et / ou publication posté dans les 12 mois"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value2 = $map[$label]
}
